$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; existing rows 3-9 shift down to 4-10.
$ws.Rows(3).Insert()

# Fill the newly inserted row 3 with the new weekly record.
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("D3").Value = 44687
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112013
$ws.Range("G3").Value = "Alcachofa"
$ws.Range("H3").Value = "Española"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 19000
$ws.Range("M3").Value = 18500
$ws.Range("N3").Value = "$/caja 30 unidades"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 617
$ws.Range("Q3").Value = 30
$ws.Range("R3").Value = "Hortaliza"
